$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.493.46"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.78"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.82%  "
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.57"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4538"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07116"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8759"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07747"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.36"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.831.45"
$ws.Range("E13").Value = "  -5.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.266"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.19"
$ws.Range("E16").Value = "  -5.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008566"
$ws.Range("E18").Value = "  -3.94%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.533.10"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.971"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.39"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.974"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.65"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.90"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.995"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.48"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.840"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08655"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.040"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7279"
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.431"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.112"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.554"
$ws.Range("E36").Value = "  -6.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.079"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01929"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05093"
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.879"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.926"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4992"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1566"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.115"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4596"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.65"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.923"
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.585"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05989"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.80"
$ws.Range("E51").Value = "  -2.85%  "
